# auto: removing some labels from the patient card
# Removes the "nick" (Nickname) and "gender_n" (Gender Identity) note rows
# from the survey sheet. Deleting the entire rows shifts everything below
# up automatically (including the specially-styled
# "n_special_instruction_title" row and its data-validation range).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Delete the lower row first so the "nick" row reference (49) below isn't
# shifted by this delete.
$ws.Rows("51").EntireRow.Delete()   # gender_n / "Gender Identity: **${gender_ctx}**"
$ws.Rows("49").EntireRow.Delete()   # nick / "Nickname: **${aka_ctx}**"
